$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 27: 221127_cat_v2data / mangan_steakNfries ---
# --- Row 28: 221127_cat_v2data_deeper / trial_JM_MS ---

# Dates: copy format+value from the last existing data row (26) which already
# carries the correct date style (numFmtId 14), then overwrite with the new date.
$ws.Range("A26").Copy()
$ws.Range("A27").PasteSpecial(-4122)
$ws.Range("A27").Value = 44892

$ws.Range("A26").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("A28").Value = 44892

# Name Ramp / Name columns - plain text values (order matches the source order
# new submissions were recorded in: Name Ramp for row 28, Name Ramp for row 27,
# Name for row 28, Name for row 27).
$ws.Range("C28").Value = "trial_JM_MS"
$ws.Range("C27").Value = "mangan_steakNfries"
$ws.Range("B28").Value = "221127_cat_v2data_deeper"
$ws.Range("B27").Value = "221127_cat_v2data"

# Hand in / By columns - identical to the prior row, so copy value+format
# straight across (keeps them as text "TRUE"/"Maria" sharing the same string).
$ws.Range("D26").Copy()
$ws.Range("D27").PasteSpecial(-4104)
$ws.Range("E26").Copy()
$ws.Range("E27").PasteSpecial(-4104)

$ws.Range("D26").Copy()
$ws.Range("D28").PasteSpecial(-4104)
$ws.Range("E26").Copy()
$ws.Range("E28").PasteSpecial(-4104)

# Extend the table and sheet dimension to include the two new rows.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E28"))

# Update the view so the new rows are visible/selected, matching the saved state.
[void]$ws.Range("B28").Select()
$ws.Application.ActiveWindow.ScrollRow = 16
